$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the sheet and scan column G for the
# "Recorded By" value "System, dnasr281@gmail.com", swapping it to
# "dnasr281@gmail.com, System" (order of the two names reversed).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
